$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21, pushing current rows 21-24 down to 22-25
$ws.Rows.Item(21).Insert()

# Populate the new row 21 with data (same as old row21 except D and M)
$ws.Range("A21").Value = 6
$ws.Range("B21").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C21").Value = "Metropolitana"
$ws.Range("D21").Value = 44641
$ws.Range("D21").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E21").Value = 13
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100101
$ws.Range("H21").Value = "Berries"
$ws.Range("I21").Value = 100101006
$ws.Range("J21").Value = "Higo"
$ws.Range("K21").Value = "Sin especificar"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 13000
$ws.Range("O21").Value = 13000
$ws.Range("P21").Value = 13000
$ws.Range("Q21").Value = "$/bandeja 7 kilos"
$ws.Range("R21").Value = "Región Metropolitana"
$ws.Range("S21").Value = 1857
$ws.Range("T21").Value = 7
